$wb = $excel.ActiveWorkbook

# --- Sheet "Студенты" (Students) ---
$students = $wb.Worksheets.Item(1)

# Row 8: remove the student's name ("Витальев В. А.") entirely
$students.Range("B8").ClearContents()

# Row 9: replace the student's name ("Петров П. А.") with a blank,
# quote-prefixed (text-forced) value of five spaces
$students.Range("B9").Value = "'     "

# Update the selection on this sheet
$students.Range("B9").Select()

# --- Sheet "Университеты" (Universities) ---
$universities = $wb.Worksheets.Item(2)

# Row 4: remove the university's full name
# ("Московский Государственный Медицинский Университет") entirely
$universities.Range("B4").ClearContents()

# Row 5: replace the university's abbreviation ("ТУМ") with a blank,
# quote-prefixed (text-forced) value of five spaces
$universities.Range("C5").Value = "'     "

# Make this sheet the active / selected tab with the given selection
$universities.Activate()
$universities.Range("C6").Select()
